$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.861.06"
$ws.Range("E2").Value = "  +1.01%  "
$ws.Range("D3").Value = "2.285.44"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.91"
$ws.Range("E5").Value = "  -3.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.95"
$ws.Range("E6").Value = "  -1.50%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.616"
$ws.Range("E7").Value = "  -1.09%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.602"
$ws.Range("E9").Value = "  -0.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.81"
$ws.Range("E10").Value = "  -3.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0898"
$ws.Range("E11").Value = "  -1.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.21"
$ws.Range("E12").Value = "  -3.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.108"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.980"
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.25"
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("D16").Value = "2.629.67"
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("D17").Value = "2.286.29"
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("D18").Value = "42.520.56"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.29"
$ws.Range("E19").Value = "  -2.57%  "
$ws.Range("E20").Value = "  -1.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.44"
$ws.Range("E21").Value = "  +0.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.22"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "268.26"
$ws.Range("E23").Value = "  -0.49%  "
$ws.Range("E24").Value = "  -5.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.21"
$ws.Range("E25").Value = "  -0.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.78"
$ws.Range("E27").Value = "  -0.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.13"
$ws.Range("E28").Value = "  +15.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.30"
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.38"
$ws.Range("E30").Value = "  -0.82%  "
$ws.Range("E31").Value = "  -6.95%  "
$ws.Range("E32").Value = "  -0.72%  "
$ws.Range("E33").Value = "  -3.79%  "
$ws.Range("E34").Value = "  -1.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.56"
$ws.Range("E35").Value = "  +1.07%  "
$ws.Range("E36").Value = "  -2.42%  "
$ws.Range("E37").Value = "  -1.91%  "
$ws.Range("E38").Value = "  -2.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.77"
$ws.Range("E39").Value = "  +1.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.61"
$ws.Range("E40").Value = "  -3.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "107.56"
$ws.Range("E41").Value = "  +11.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.55"
$ws.Range("E42").Value = "  +0.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.91"
$ws.Range("E43").Value = "  +1.44%  "
$ws.Range("E44").Value = "  +0.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  -0.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.03"
$ws.Range("E46").Value = "  -3.06%  "
$ws.Range("D47").Value = "1.727.98"
$ws.Range("E47").Value = "  +8.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "110.29"
$ws.Range("E48").Value = "  -2.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "76.43"
$ws.Range("E49").Value = "  -6.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.64"
$ws.Range("E50").Value = "  -2.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.15"
$ws.Range("E51").Value = "  -2.28%  "
